$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.524.50"
$ws.Range("E2").Value = "  -0.25%  "

$ws.Range("D3").Value = "1.811.46"
$ws.Range("E3").Value = "  +0.54%  "

$ws.Range("E4").Value = "  +0.24%  "

$ws.Range("D5").Value = "'228.75"
$ws.Range("E5").Value = "  +0.52%  "

$ws.Range("D6").Value = "'0.578"
$ws.Range("E6").Value = "  +3.81%  "

$ws.Range("E7").Value = "  +0.20%  "

$ws.Range("D8").Value = "'34.93"
$ws.Range("E8").Value = "  +6.66%  "

$ws.Range("E9").Value = "  +1.50%  "

$ws.Range("E10").Value = "  -0.21%  "

$ws.Range("D11").Value = "'0.0956"
$ws.Range("E11").Value = "  +0.71%  "

$ws.Range("E12").Value = "  +0.61%  "

$ws.Range("D13").Value = "'11.22"
$ws.Range("E13").Value = "  +0.40%  "

$ws.Range("D14").Value = "1.816.83"
$ws.Range("E14").Value = "  +1.00%  "

$ws.Range("D15").Value = "'0.649"
$ws.Range("E15").Value = "  +1.19%  "

$ws.Range("D16").Value = "34.532.21"
$ws.Range("E16").Value = "  -0.15%  "

$ws.Range("D17").Value = "'4.45"
$ws.Range("E17").Value = "  +2.67%  "

$ws.Range("D18").Value = "'69.22"
$ws.Range("E18").Value = "  +0.45%  "

$ws.Range("E19").Value = "  -0.87%  "

$ws.Range("D20").Value = "'245.23"
$ws.Range("E20").Value = "  -0.95%  "

$ws.Range("D21").Value = "'11.46"
$ws.Range("E21").Value = "  +1.00%  "

$ws.Range("E22").Value = "  +0.17%  "

$ws.Range("E23").Value = "  -0.76%  "

$ws.Range("D24").Value = "'172.59"
$ws.Range("E24").Value = "  +1.70%  "

$ws.Range("E25").Value = "  +1.84%  "

$ws.Range("D26").Value = "'7.99"
$ws.Range("E26").Value = "  +9.33%  "

$ws.Range("D27").Value = "'16.82"
$ws.Range("E27").Value = "  +1.25%  "

$ws.Range("D28").Value = "'0.119"
$ws.Range("E28").Value = "  +2.24%  "

$ws.Range("E29").Value = "  +0.03%  "

$ws.Range("D30").Value = "'4.03"
$ws.Range("E30").Value = "  -2.62%  "

$ws.Range("D31").Value = "'0.0533"
$ws.Range("E31").Value = "  +1.08%  "

$ws.Range("D32").Value = "'3.86"
$ws.Range("E32").Value = "  +0.71%  "

$ws.Range("E33").Value = "  +0.01%  "

$ws.Range("D34").Value = "'1.84"
$ws.Range("E34").Value = "  -0.12%  "

$ws.Range("D35").Value = "'0.683"
$ws.Range("E35").Value = "  +0.57%  "

$ws.Range("D36").Value = "1.395.70"
$ws.Range("E36").Value = "  -2.55%  "

$ws.Range("E37").Value = "  -5.59%  "

$ws.Range("E38").Value = "  -1.55%  "

$ws.Range("E39").Value = "  -0.56%  "

$ws.Range("D40").Value = "'83.66"
$ws.Range("E40").Value = "  -1.99%  "

$ws.Range("D41").Value = "'0.962"
$ws.Range("E41").Value = "  +1.34%  "

$ws.Range("D42").Value = "'2.83"
$ws.Range("E42").Value = "  +2.47%  "

$ws.Range("E43").Value = "  -0.82%  "

$ws.Range("D44").Value = "'13.39"
$ws.Range("E44").Value = "  -3.37%  "

$ws.Range("E45").Value = "  +3.93%  "

$ws.Range("D46").Value = "'0.0515"
$ws.Range("E46").Value = "  -2.01%  "

$ws.Range("E47").Value = "  -1.78%  "

$ws.Range("D48").Value = "1.973.65"
$ws.Range("E48").Value = "  +0.68%  "

$ws.Range("D49").Value = "'105.19"
$ws.Range("E49").Value = "  -0.68%  "

$ws.Range("E50").Value = "  +2.64%  "
